# "Se limpio la plantilla" - clear stray "SÍ" answers that were left in the
# checklist template so it ships blank, and leave the workbook focused on the
# "Resumen" (summary) sheet instead of "Productos".

$wb = $excel.ActiveWorkbook

# Clear the leftover "SÍ" marks in the "procesos" checklist (rows 13-17).
$wsProcesos = $wb.Worksheets.Item("procesos")
$wsProcesos.Range("C13:C17").ClearContents()

# Clear the leftover "SÍ" marks in the "Productos" checklist (rows 37-45).
$wsProductos = $wb.Worksheets.Item("Productos")
$wsProductos.Range("C37:C45").ClearContents()

# Reset each checklist sheet's selection back to the top summary row.
$wsProcesos.Range("C1:E1").Select()
$wsProductos.Range("C1:E1").Select()

# Make "Resumen" the active/selected sheet again (it had lost focus to
# "Productos"), and park the selection back on its last-used cell.
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Activate()
$wsResumen.Range("D43").Select()
